# Update '想去人数' (F column) values per the target diff.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("展览")
$ws.Cells.Item(3, 6).Value = 164
$ws.Cells.Item(4, 6).Value = 42
$ws.Cells.Item(6, 6).Value = 5863
$ws.Cells.Item(8, 6).Value = 432
$ws.Cells.Item(9, 6).Value = 3957
$ws.Cells.Item(10, 6).Value = 81
$ws.Cells.Item(11, 6).Value = 28
$ws.Cells.Item(17, 6).Value = 134
$ws.Cells.Item(18, 6).Value = 660
$ws.Cells.Item(19, 6).Value = 3960
$ws.Cells.Item(22, 6).Value = 5488
$ws.Cells.Item(24, 6).Value = 2164
$ws.Cells.Item(25, 6).Value = 140
$ws.Cells.Item(26, 6).Value = 377
$ws.Cells.Item(27, 6).Value = 8184
$ws.Cells.Item(29, 6).Value = 2222
$ws.Cells.Item(30, 6).Value = 2246
$ws.Cells.Item(32, 6).Value = 182
$ws.Cells.Item(33, 6).Value = 1342
$ws.Cells.Item(37, 6).Value = 261
$ws.Cells.Item(38, 6).Value = 24
$ws.Cells.Item(39, 6).Value = 16
$ws.Cells.Item(40, 6).Value = 1194
$ws.Cells.Item(41, 6).Value = 1188
$ws.Cells.Item(46, 6).Value = 2182
$ws.Cells.Item(48, 6).Value = 239

$ws = $wb.Worksheets.Item("演出")
$ws.Cells.Item(6, 6).Value = 154
$ws.Cells.Item(19, 6).Value = 10
$ws.Cells.Item(22, 6).Value = 20

$ws = $wb.Worksheets.Item("本地生活")
$ws.Cells.Item(2, 6).Value = 607
$ws.Cells.Item(3, 6).Value = 788
$ws.Cells.Item(4, 6).Value = 74

$ws = $wb.Worksheets.Item("全部类型")
$ws.Cells.Item(3, 6).Value = 164
$ws.Cells.Item(4, 6).Value = 42
$ws.Cells.Item(5, 6).Value = 607
$ws.Cells.Item(6, 6).Value = 788
$ws.Cells.Item(7, 6).Value = 5863
$ws.Cells.Item(8, 6).Value = 432
$ws.Cells.Item(9, 6).Value = 3957
$ws.Cells.Item(10, 6).Value = 81
$ws.Cells.Item(11, 6).Value = 28
$ws.Cells.Item(16, 6).Value = 154
$ws.Cells.Item(18, 6).Value = 660
$ws.Cells.Item(19, 6).Value = 3960
$ws.Cells.Item(23, 6).Value = 5489
$ws.Cells.Item(25, 6).Value = 2164
$ws.Cells.Item(26, 6).Value = 140
$ws.Cells.Item(27, 6).Value = 377
$ws.Cells.Item(28, 6).Value = 8184
$ws.Cells.Item(30, 6).Value = 2222
$ws.Cells.Item(31, 6).Value = 2246
$ws.Cells.Item(32, 6).Value = 182
$ws.Cells.Item(33, 6).Value = 1342
$ws.Cells.Item(36, 6).Value = 261
$ws.Cells.Item(37, 6).Value = 16
$ws.Cells.Item(38, 6).Value = 1194
$ws.Cells.Item(39, 6).Value = 1188
$ws.Cells.Item(44, 6).Value = 2182
$ws.Cells.Item(47, 6).Value = 239
$ws.Cells.Item(48, 6).Value = 20
